# ADD results from server
# Update investment-cost figures in row 2 of the "2025", "2030", "2035",
# "2040" and "2045" sheets. ("2050" sheet is unchanged.)

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1812.803126660054
$ws.Range("E2").Value = 4622.8361210227
$ws.Range("G2").Value = 2428.77771379855
$ws.Range("I2").Value = 10371.325902912
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1632.534049630441
$ws.Range("O2").Value = 2420.79197950242

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 1358.894601537055
$ws.Range("B2").Value = 3025.951641828148
$ws.Range("E2").Value = 12889.56086624117
$ws.Range("I2").Value = 19330.19497335958
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2773.846293838942
$ws.Range("O2").Value = 2340.403967713372

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 1770.119671271456
$ws.Range("B2").Value = 4686.043832468928
$ws.Range("E2").Value = 19385.71650516988
$ws.Range("I2").Value = 28665.83942972087
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 5043.434707998173
$ws.Range("O2").Value = 3834.906109154087

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("O2").Value = 247.0922343210123

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 775.978146433487
$ws.Range("O2").Value = 204.0473635070248

# "2050" sheet is unchanged.
